# The commit swaps the contents of ppt/theme/theme1.xml ("Office Theme" /
# "Office" colour scheme) and ppt/theme/theme2.xml ("Integral" / "Red Violet"
# colour scheme), while every relationship keeps pointing at the same target
# filenames it always did (theme2.xml stays the slide master's / the
# presentation's theme part, theme1.xml stays the notes master's theme part).
#
# theme1.xml and theme2.xml are byte-for-byte identical once the <a:theme>
# name, the <a:clrScheme> name and the twelve colour values are normalised
# away (fontScheme / fmtScheme are untouched by the edit), so the swap is
# reproduced here by rewriting the twelve colours that PowerPoint's object
# model exposes through Slide.ThemeColorScheme (it is a single shared part,
# so editing it from any slide updates ppt/theme/theme2.xml for the whole
# deck) to the values theme2.xml's counterpart ("Office") used to hold.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index order matches <a:clrScheme>: dk1, lt1, dk2, lt2, accent1..accent6,
# hlink, folHlink. Values are the "Office" scheme's RGB()-style integers
# (R + G*256 + B*65536).
$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
